$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new column G to the existing header + data rows (1-4) ---
$ws.Cells.Item(1,7).Value = "xTestWord2x"
$ws.Cells.Item(2,7).Value = "Yowz"
$ws.Cells.Item(3,7).Value = "Hahahaha"
$ws.Cells.Item(4,7).Value = "Wassap"

# --- Duplicate rows 2-4 (now A:G) into rows 5-7, carrying the existing ---
# --- cell-level formatting down (so new cells inherit style index 1   ---
# --- just like the source rows do) instead of picking up the default. ---
$ws.Range("A2:G4").Copy()
$ws.Range("A5:G7").Insert(-4121)

# --- Fix up the first column of the new rows with the new document names ---
$ws.Cells.Item(5,1).Value = "Test Doc 4"
$ws.Cells.Item(6,1).Value = "Test Doc 5"
$ws.Cells.Item(7,1).Value = "Test Doc 6"

# --- Column G of the new rows mirrors column F of the "paired" row, ---
# --- exactly like rows 2-4 do. ---
$ws.Cells.Item(5,7).Value = "Yowz"
$ws.Cells.Item(6,7).Value = "Hahahaha"
$ws.Cells.Item(7,7).Value = "Wassap"

# --- Give the new column F a fixed, custom width like the other data columns ---
$ws.Columns(6).ColumnWidth = 11.666

# --- Update the active selection to mirror the saved view state ---
$ws.Range("H11").Select() | Out-Null
